# Apply the "pushing new excel sheet" edit:
#  1. Append two new worksheets: "Invalidcodeveni" and "tryeditorcode"
#  2. Rewrite the "Login" sheet test-data table (new header names, new/changed rows)
#  3. Make "Login" the active/selected sheet & tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$invalidVeni = $wb.Worksheets.Add($null, $lastSheet)
$invalidVeni.Name = "Invalidcodeveni"
$invalidVeni.Range("A1").Value = "pythonCode"
$invalidVeni.Range("B1").Value = "output"
$invalidVeni.Range("A2").Value = "hello"
$invalidVeni.Range("B2").Value = "NameError: name 'hello' is not defined on line 1"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tryEditor = $wb.Worksheets.Add($null, $lastSheet)
$tryEditor.Name = "tryeditorcode"
$tryEditor.Range("A1").Value = "pythonCode"
$tryEditor.Range("B1").Value = "output"
$tryEditor.Range("A2").Value = 'print("Hello")'
$tryEditor.Range("B2").Value = "Hello"
$tryEditor.Range("A3").Value = "Hello"

# ---------------------------------------------------------------------------
# 2. Rewrite the "Login" sheet
# ---------------------------------------------------------------------------
$login = $wb.Worksheets.Item("Login")

# Clear the previous table (old range was A1:C5) before writing the new one.
$login.Range("A1:C11").ClearContents()

# Header row
$login.Range("A1").Value = "Username"
$login.Range("B1").Value = "password"
$login.Range("C1").Value = "Expected Message"

# Data rows
$login.Range("A2").Value = "Rockstars_Numpy"
$login.Range("B2").Value = "Numpy@Rock123"
$login.Range("C2").Value = "You are logged in"

$login.Range("A3").Value = "Rockstars_Numpy"
$login.Range("C3").Value = "Please fill out this field."

$login.Range("B4").Value = "Numpy@Rock123"
$login.Range("C4").Value = "Please fill out this field."

$login.Range("C5").Value = "Please fill out this field."

$login.Range("A6").Value = "username"
$login.Range("B6").Value = "Numpy@Rock123"
$login.Range("C6").Value = "Invalid Username and Password"

$login.Range("A7").Value = "Rockstars_Numpy"
$login.Range("B7").Value = "sdet84batch"
$login.Range("C7").Value = "Invalid Username and Password"

# Match the original look-and-feel: the sheet alternates between a black
# Arial-14 font (style "7") and a dark-grey Arial-14 font (style "8").
$login.Range("A1:C7").Font.Name = "Arial"
$login.Range("A1:C7").Font.Size = 14
$login.Range("A1:C7").Font.Color = 0

$greyCells = "A2", "B2", "A3", "B3", "B6", "A7"
foreach ($ref in $greyCells) {
    $login.Range($ref).Font.Color = 3355443
}

# ---------------------------------------------------------------------------
# 3. Make "Login" the active tab / selected sheet
# ---------------------------------------------------------------------------
$login.Activate()
$login.Range("G7").Select()
